$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 15 ---
$ws.Cells.Item(15, 1).Value = 112481511
$ws.Cells.Item(15, 2).Value = 90813
$ws.Cells.Item(15, 3).Value = "Ovaliderad"
$ws.Cells.Item(15, 4).Value = "VU"
$ws.Cells.Item(15, 5).Value = 1435
$ws.Cells.Item(15, 6).Value = "Bitter taggsvamp"
$ws.Cells.Item(15, 7).Value = "Hydnellum fennicum"
$ws.Cells.Item(15, 8).Value = "(P.Karst.) E.Larss., K.H.Larss. & Kõljalg"
$ws.Cells.Item(15, 9).NumberFormat = "@"
$ws.Cells.Item(15, 9).Value = "1"
$ws.Cells.Item(15, 9).Style = "Normal"
$ws.Cells.Item(15, 10).Value = "fruktkroppar"
$ws.Cells.Item(15, 16).Value = "Barrtjärnen, Bergsjövägen, Ede, Nordanstigs kommun (Barrtjärnen, Bergsjövägen, Ede, Nordanstigs kommun), Hls"
$ws.Cells.Item(15, 17).Value = 601183
$ws.Cells.Item(15, 18).Value = 6877672
$ws.Cells.Item(15, 19).Value = 25
$ws.Cells.Item(15, 20).Value = "Gävleborg"
$ws.Cells.Item(15, 21).Value = "Nordanstig"
$ws.Cells.Item(15, 22).Value = "Hälsingland"
$ws.Cells.Item(15, 23).Value = "Bergsjö"
$ws.Cells.Item(15, 25).NumberFormat = "@"
$ws.Cells.Item(15, 25).Value = "2023-09-09"
$ws.Cells.Item(15, 25).Style = "Normal"
$ws.Cells.Item(15, 26).Value = "15:00"
$ws.Cells.Item(15, 27).NumberFormat = "@"
$ws.Cells.Item(15, 27).Value = "2023-09-09"
$ws.Cells.Item(15, 27).Style = "Normal"
$ws.Cells.Item(15, 28).Value = "15:00"
$ws.Cells.Item(15, 29).Value = "I slänten ner mot en större svacka i terrängen. Ca 3 m från stigen. Tydlig doft av bittermandel."
$ws.Cells.Item(15, 30).Value = $false
$ws.Cells.Item(15, 31).Value = $false
$ws.Cells.Item(15, 33).Value = $false
$ws.Cells.Item(15, 49).Value = "Henrik Tykosson"
$ws.Cells.Item(15, 50).Value = "Henrik Tykosson"

# --- Row 16 ---
$ws.Cells.Item(16, 1).Value = 112481246
$ws.Cells.Item(16, 2).Value = 90291
$ws.Cells.Item(16, 3).Value = "Ovaliderad"
$ws.Cells.Item(16, 4).Value = "VU"
$ws.Cells.Item(16, 5).Value = 1958
$ws.Cells.Item(16, 6).Value = "Lammticka"
$ws.Cells.Item(16, 7).Value = "Albatrellus subrubescens"
$ws.Cells.Item(16, 8).Value = "(Murrill) Pouzar"
$ws.Cells.Item(16, 16).Value = "Barrtjärnen, Bergsjövägen, Ede, Nordanstigs kommun (Barrtjärnen, Bergsjövägen, Ede, Nordanstigs kommun), Hls"
$ws.Cells.Item(16, 17).Value = 601210
$ws.Cells.Item(16, 18).Value = 6877659
$ws.Cells.Item(16, 19).Value = 25
$ws.Cells.Item(16, 20).Value = "Gävleborg"
$ws.Cells.Item(16, 21).Value = "Nordanstig"
$ws.Cells.Item(16, 22).Value = "Hälsingland"
$ws.Cells.Item(16, 23).Value = "Bergsjö"
$ws.Cells.Item(16, 25).NumberFormat = "@"
$ws.Cells.Item(16, 25).Value = "2023-10-02"
$ws.Cells.Item(16, 25).Style = "Normal"
$ws.Cells.Item(16, 26).Value = "16:05"
$ws.Cells.Item(16, 27).NumberFormat = "@"
$ws.Cells.Item(16, 27).Value = "2023-10-02"
$ws.Cells.Item(16, 27).Style = "Normal"
$ws.Cells.Item(16, 28).Value = "16:05"
$ws.Cells.Item(16, 29).Value = "I kanten av stigen ned mot båtplats. Västra kanten I böjen."
$ws.Cells.Item(16, 30).Value = $false
$ws.Cells.Item(16, 31).Value = $false
$ws.Cells.Item(16, 33).Value = $false
$ws.Cells.Item(16, 49).Value = "Henrik Tykosson"
$ws.Cells.Item(16, 50).Value = "Henrik Tykosson"
